$p = $ppt.ActivePresentation

# Insert a new "Title and Content" slide before the current slide 7
# ("Practice - Assignment 1-2"), pushing it and the following slide down.
$s = $p.Slides.Add(7, 2)

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Video Time!"

# Content placeholder: hyperlinked YouTube URL followed by a blank line
$url = "https://www.youtube.com/watch?v=u05YFFh6044"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$tr.Text = $url + "`r"
$link = $tr.Characters(1, $url.Length)
$link.ActionSettings.Item(1).Hyperlink.Address = $url
